$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2830598511410756
$ws.Range("C2").Value = 0.06452926745294008
$ws.Range("E2").Value = 0.421958101597582
$ws.Range("F2").Value = 0.4443680307746121
$ws.Range("G2").Value = 0.2614843200896999
$ws.Range("H2").Value = 0.4448408548528633
$ws.Range("I2").Value = 0.3275972041486064
$ws.Range("K2").Value = 0.3034573273221213
$ws.Range("N2").Value = 1.024950686546191
$ws.Range("O2").Value = 1.332988765104815

$ws.Range("B3").Value = 0.2472426830551626
$ws.Range("C3").Value = 0.05941881472426758
$ws.Range("E3").Value = 0.3682353146068067
$ws.Range("F3").Value = 0.3878228170618172
$ws.Range("G3").Value = 0.2631787881475347
$ws.Range("H3").Value = 0.4488222816131326
$ws.Range("I3").Value = 0.3320083991610137
$ws.Range("K3").Value = 0.2649728001601943
$ws.Range("N3").Value = 1.033539013173417
$ws.Range("O3").Value = 1.344704553376033

$ws.Range("B4").Value = 0.2251901798649101
$ws.Range("C4").Value = 0.05625834988970269
$ws.Range("E4").Value = 0.3353193148865472
$ws.Range("F4").Value = 0.3531389305169483
$ws.Range("G4").Value = 0.2644492813688899
$ws.Range("H4").Value = 0.4514755866007647
$ws.Range("I4").Value = 0.3349209906004429
$ws.Range("K4").Value = 0.241246678301593
$ws.Range("N4").Value = 1.039184287165835
$ws.Range("O4").Value = 1.352815633331844

$ws.Range("B5").Value = 0.2161889182791867
$ws.Range("C5").Value = 0.05496479153846678
$ws.Range("E5").Value = 0.3219221104235714
$ws.Range("F5").Value = 0.3390132514313251
$ws.Range("G5").Value = 0.2650247193501727
$ws.Range("H5").Value = 0.4526093027551141
$ws.Range("I5").Value = 0.3361591594807223
$ws.Range("K5").Value = 0.2315543691354804
$ws.Range("N5").Value = 1.041578402895173
$ws.Range("O5").Value = 1.356351344368463

$ws.Range("B6").Value = 0.2146933949222785
$ws.Range("C6").Value = 0.05474965768905804
$ws.Range("E6").Value = 0.3196984642130616
$ws.Range("F6").Value = 0.336668177824194
$ws.Range("G6").Value = 0.2651237509034061
$ws.Range("H6").Value = 0.4528007245607526
$ws.Range("I6").Value = 0.3363678510065142
$ws.Range("K6").Value = 0.2299435516247001
$ws.Range("N6").Value = 1.041981600751896
$ws.Range("O6").Value = 1.356952350877833

$ws.Range("B7").Value = 0.2250688444366347
$ws.Range("C7").Value = 0.0562409272977078
$ws.Range("E7").Value = 0.3351385710433163
$ws.Range("F7").Value = 0.3529483938344953
$ws.Range("G7").Value = 0.264456808500853
$ws.Range("H7").Value = 0.4514906638393015
$ws.Range("I7").Value = 0.3349374814899235
$ws.Range("K7").Value = 0.2411160597690127
$ws.Range("N7").Value = 1.039216195879266
$ws.Range("O7").Value = 1.352862384824164

$ws.Range("B8").Value = 0.2707229947996836
$ws.Range("C8").Value = 0.06277190566703439
$ws.Range("E8").Value = 0.4034191731270909
$ws.Range("F8").Value = 0.4248636149813478
$ws.Range("G8").Value = 0.2620207255325369
$ws.Range("H8").Value = 0.4461703337591629
$ws.Range("I8").Value = 0.3290757884598641
$ws.Range("K8").Value = 0.290208193630292
$ws.Range("N8").Value = 1.027834777289677
$ws.Range("O8").Value = 1.33683772678053

$ws.Range("B9").Value = 0.3597494917518702
$ws.Range("C9").Value = 0.07539815445854003
$ws.Range("E9").Value = 0.5379362583876031
$ws.Range("F9").Value = 0.5661985755041457
$ws.Range("G9").Value = 0.259076384846864
$ws.Range("H9").Value = 0.4373933171661477
$ws.Range("I9").Value = 0.3192030159756687
$ws.Range("K9").Value = 0.385693239844727
$ws.Range("N9").Value = 1.008464164995068
$ws.Range("O9").Value = 1.312709704436614

$ws.Range("B10").Value = 0.4248308193372736
$ws.Range("C10").Value = 0.08456323805377508
$ws.Range("E10").Value = 0.6372402226762546
$ws.Range("F10").Value = 0.6702781546542269
$ws.Range("G10").Value = 0.2580408170021329
$ws.Range("H10").Value = 0.4319548719334207
$ws.Range("I10").Value = 0.3129417071646508
$ws.Range("K10").Value = 0.455348829125569
$ws.Range("N10").Value = 0.9960253494081215
$ws.Range("O10").Value = 1.299453876779353

$ws.Range("B11").Value = 0.4543630593797729
$ws.Range("C11").Value = 0.08870830162301502
$ws.Range("E11").Value = 0.6825405946669605
$ws.Range("F11").Value = 0.7176906081379002
$ws.Range("G11").Value = 0.2578168075535743
$ws.Range("H11").Value = 0.4297002065642488
$ws.Range("I11").Value = 0.3103094616895987
$ws.Range("K11").Value = 0.4869253491602308
$ws.Range("N11").Value = 0.9907549757332319
$ws.Range("O11").Value = 1.294399339933463

$ws.Range("B12").Value = 0.4655350867656125
$ws.Range("C12").Value = 0.09027441789704937
$ws.Range("E12").Value = 0.6997144006774505
$ws.Range("F12").Value = 0.7356546913071611
$ws.Range("G12").Value = 0.2577676866868899
$ws.Range("H12").Value = 0.4288779725042744
$ws.Range("I12").Value = 0.3093438334091729
$ws.Range("K12").Value = 0.4988662589130684
$ws.Range("N12").Value = 0.9888149717105179
$ws.Range("O12").Value = 1.292626014630073

$ws.Range("B13").Value = 0.4631294981970484
$ws.Range("C13").Value = 0.08993728453690153
$ws.Range("E13").Value = 0.6960148237615016
$ws.Range("F13").Value = 0.7317853510981394
$ws.Range("G13").Value = 0.2577766748561316
$ws.Range("H13").Value = 0.4290536513773517
$ws.Range("I13").Value = 0.3095504122768205
$ws.Range("K13").Value = 0.4962953106370662
$ws.Range("N13").Value = 0.9892303063681496
$ws.Range("O13").Value = 1.293001666676759

$ws.Range("B14").Value = 0.4552824165704692
$ws.Range("C14").Value = 0.0888372179826149
$ws.Range("E14").Value = 0.683953094698623
$ws.Range("F14").Value = 0.7191683204515869
$ws.Range("G14").Value = 0.2578120497196252
$ws.Range("H14").Value = 0.4296319282599015
$ws.Range("I14").Value = 0.3102293942142751
$ws.Range("K14").Value = 0.4879080677260959
$ws.Range("N14").Value = 0.9905942528810101
$ws.Range("O14").Value = 1.294250624743938

$ws.Range("B15").Value = 0.4504743770137623
$ws.Range("C15").Value = 0.08816293427967992
$ws.Range("E15").Value = 0.6765675219679395
$ws.Range("F15").Value = 0.7114413442032514
$ws.Range("G15").Value = 0.2578383731709764
$ws.Range("H15").Value = 0.4299902500981432
$ws.Range("I15").Value = 0.3106493486802329
$ws.Range("K15").Value = 0.4827684816247029
$ws.Range("N15").Value = 0.991436971626527
$ws.Range("O15").Value = 1.295033985964196

$ws.Range("B16").Value = 0.4228992851685973
$ws.Range("C16").Value = 0.08429185737027467
$ws.Range("E16").Value = 0.6342823973682101
$ws.Range("F16").Value = 0.6671810134426437
$ws.Range("G16").Value = 0.2580604435123277
$ws.Range("H16").Value = 0.4321066345942555
$ws.Range("I16").Value = 0.313118084800486
$ws.Range("K16").Value = 0.4532829619531071
$ws.Range("N16").Value = 0.996377587173285
$ws.Range("O16").Value = 1.299803870830971

$ws.Range("B17").Value = 0.4059635789577953
$ws.Range("C17").Value = 0.08191084686997385
$ws.Range("E17").Value = 0.6083751992738655
$ws.Range("F17").Value = 0.6400460337125793
$ws.Range("G17").Value = 0.2582600847636556
$ws.Range("H17").Value = 0.4334611542938589
$ws.Range("I17").Value = 0.3146879664658009
$ws.Range("K17").Value = 0.4351659147870635
$ws.Range("N17").Value = 0.9995078707544351
$ws.Range("O17").Value = 1.302980220451758

$ws.Range("B18").Value = 0.3962157191526217
$ws.Range("C18").Value = 0.08053907988750098
$ws.Range("E18").Value = 0.5934858899798314
$ws.Range("F18").Value = 0.6244449056556647
$ws.Range("G18").Value = 0.258398158444578
$ws.Range("H18").Value = 0.4342608793967599
$ws.Range("I18").Value = 0.3156112523830288
$ws.Range("K18").Value = 0.4247351352995281
$ws.Range("N18").Value = 1.001344859470123
$ws.Range("O18").Value = 1.304898975765482

$ws.Range("B19").Value = 0.3929140966349962
$ws.Range("C19").Value = 0.08007423441850392
$ws.Range("E19").Value = 0.5884466158371424
$ws.Range("F19").Value = 0.619163680173358
$ws.Range("G19").Value = 0.2584488948169792
$ws.Range("H19").Value = 0.4345351970544655
$ws.Range("I19").Value = 0.3159273507292788
$ws.Range("K19").Value = 0.4212016959421305
$ws.Range("N19").Value = 1.001973107937623
$ws.Range("O19").Value = 1.305564387604591

$ws.Range("B20").Value = 0.4077671312568896
$ws.Range("C20").Value = 0.08216454510075266
$ws.Range("E20").Value = 0.6111318310836538
$ws.Range("F20").Value = 0.642933953830422
$ws.Range("G20").Value = 0.2582364255222984
$ws.Range("H20").Value = 0.4333148269122162
$ws.Range("I20").Value = 0.3145187450866018
$ws.Range("K20").Value = 0.4370955791240192
$ws.Range("N20").Value = 0.99917086603336
$ws.Range("O20").Value = 1.302632588247491

$ws.Range("B21").Value = 0.4575876032357655
$ws.Range("C21").Value = 0.08916043039720023
$ws.Range("E21").Value = 0.6874953758624542
$ws.Range("F21").Value = 0.7228739723491628
$ws.Range("G21").Value = 0.257800688782325
$ws.Range("H21").Value = 0.4294612176468604
$ws.Range("I21").Value = 0.3100291148501526
$ws.Range("K21").Value = 0.4903720549853006
$ws.Range("N21").Value = 0.9901921153696094
$ws.Range("O21").Value = 1.293879953027144

$ws.Range("B22").Value = 0.4900826497917024
$ws.Range("C22").Value = 0.09371205003698435
$ws.Range("E22").Value = 0.7375178820264097
$ws.Range("F22").Value = 0.7751780083420101
$ws.Range("G22").Value = 0.2577241176244698
$ws.Range("H22").Value = 0.4271266205647422
$ws.Range("I22").Value = 0.3072764604619422
$ws.Range("K22").Value = 0.5250952529825383
$ws.Range("N22").Value = 0.984649050089665
$ws.Range("O22").Value = 1.288980011487695

$ws.Range("B23").Value = 0.4727456274813164
$ws.Range("C23").Value = 0.09128466674539482
$ws.Range("E23").Value = 0.7108090035271744
$ws.Range("F23").Value = 0.7472568307830727
$ws.Range("G23").Value = 0.2577458754264015
$ws.Range("H23").Value = 0.4283557992321363
$ws.Range("I23").Value = 0.3087289638361419
$ws.Range("K23").Value = 0.5065718136214059
$ws.Range("N23").Value = 0.987577756179995
$ws.Range("O23").Value = 1.291519993462771

$ws.Range("B24").Value = 0.406951780403972
$ws.Range("C24").Value = 0.08204985713648227
$ws.Range("E24").Value = 0.6098855421311811
$ws.Range("F24").Value = 0.6416283278902171
$ws.Range("G24").Value = 0.2582470493083875
$ws.Range("H24").Value = 0.4333809161268292
$ws.Range("I24").Value = 0.3145951854730384
$ws.Range("K24").Value = 0.4362232247056568
$ws.Range("N24").Value = 0.9993231095286887
$ws.Range("O24").Value = 1.302789464258268

$ws.Range("B25").Value = 0.3357212845599236
$ws.Range("C25").Value = 0.07200190853740196
$ws.Range("E25").Value = 0.5014695623254113
$ws.Range("F25").Value = 0.5279251897347166
$ws.Range("G25").Value = 0.2596756265226148
$ws.Range("H25").Value = 0.4395903921093307
$ws.Range("I25").Value = 0.3216998658515564
$ws.Range("K25").Value = 0.3599479609389959
$ws.Range("N25").Value = 1.013389278294184
$ws.Range("O25").Value = 1.318453480784768
